# ES2N - Requisitos Nao Funcionais: fix req nao funcionais
#
# 1) Date field "22/08/2023" -> "16/09/2023", typed/edited as four
#    separate runs ("16", "/0", "9", "/2023") that keep the original
#    bold/underline Arial-28 formatting.
# 2) "DropTable" + "____" (separated by spell-check proofErr markers)
#    collapse into a single "DropTable____" run.
# 3) Header: "Disciplina: ... prof.ª " + "Denilce" + " Veloso" (split by
#    spell-check proofErr markers around "Denilce") collapse into one run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Data: 22/08/2023  ->  16/09/2023
# ---------------------------------------------------------------------
$bodyText = $d.Content.Text
$dateIdx = $bodyText.IndexOf("22/08/2023")

# "22" -> "16"
$seg1 = $d.Range($dateIdx + 0, $dateIdx + 2)
$seg1.Text = "16"
$seg1.Font.Color = 255
$seg1.Font.Color = 0

# "/0" stays the same text, but becomes its own run
$seg2 = $d.Range($dateIdx + 2, $dateIdx + 4)
$seg2.Font.Color = 255
$seg2.Font.Color = 0

# "8" -> "9"
$seg3 = $d.Range($dateIdx + 4, $dateIdx + 5)
$seg3.Text = "9"
$seg3.Font.Color = 255
$seg3.Font.Color = 0

# "/2023" stays the same text, but becomes its own run
$seg4 = $d.Range($dateIdx + 5, $dateIdx + 10)
$seg4.Font.Color = 255
$seg4.Font.Color = 0

# ---------------------------------------------------------------------
# 2) DropTable ____  ->  DropTable____  (single run, proofErr removed)
#    The proofErr spellStart marker sits exactly at the start of
#    "DropTable", so a plain Find/Replace across "DropTable____" merges
#    the two runs and drops spellEnd, but leaves a dangling spellStart.
#    Widening the match one character to the left (into "Grupo: ___")
#    clears both markers, at the cost of merging across the formatting
#    boundary; we then restore the underline on just "DropTable____".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("_DropTable____", $false, $false, $false, $false, `
    $false, $true, 1, $false, "_DropTable____", 2) | Out-Null

$bodyText2 = $d.Content.Text
$dtIdx = $bodyText2.IndexOf("DropTable____")
$dtRun = $d.Range($dtIdx, $dtIdx + 13)
$dtRun.Font.Underline = 1

# ---------------------------------------------------------------------
# 3) Header: merge "Disciplina... prof.ª " + "Denilce" + " Veloso"
#    The proofErr pair wrapping "Denilce" sits strictly inside this
#    wider match, so both markers are cleanly dropped by the replace.
# ---------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$enDash = [char]0x2013
$ordFem = [char]0x00AA
$fullName = "Disciplina: Engenharia de Software 2 $enDash Turma Noite $enDash prof.$ordFem Denilce Veloso"
$hdr.Range.Find.Execute($fullName, $false, $false, $false, $false, `
    $false, $true, 1, $false, $fullName, 2) | Out-Null
